# Apply updated crypto price/volume figures to sheet1 (columns D and E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells below are prefixed with a literal apostrophe so Excel stores them as
# text (preserving leading/trailing zeros and avoiding scientific notation),
# exactly like a user typing '<value> into the cell.

$ws.Range('D2').Value = '28.681.04'
$ws.Range('E2').Value = '  -1.63%  '
$ws.Range('D3').Value = '1.803.50'
$ws.Range('D4').Value = '1.002'
$ws.Range('E4').Value = '  +0.36%  '
$ws.Range('E5').Value = '  -1.08%  '
$ws.Range('D6').Value = '0.5924'
$ws.Range('E6').Value = '  -1.50%  '
$ws.Range('E7').Value = '  +0.43%  '
$ws.Range('D8').Value = '0.2775'
$ws.Range('E8').Value = '  -0.59%  '
$ws.Range('D9').Value = '0.06822'
$ws.Range('E9').Value = '  -3.38%  '
$ws.Range('D10').Value = '23.32'
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('D11').Value = '0.07508'
$ws.Range('E11').Value = '  -1.24%  '
$ws.Range('D12').Value = '1.800.06'
$ws.Range('E12').Value = '  -1.45%  '
$ws.Range('D13').Value = '4.768'
$ws.Range('E13').Value = '  -0.26%  '
$ws.Range('D14').Value = '0.6222'
$ws.Range('E14').Value = '  -0.94%  '
$ws.Range('D15').Value = '2.048.39'
$ws.Range('E15').Value = '  -1.16%  '
$ws.Range('D16').Value = '''0.000009209'
$ws.Range('D17').Value = '''75.60'
$ws.Range('E17').Value = '  -4.08%  '
$ws.Range('D18').Value = '28.636.81'
$ws.Range('E18').Value = '  -1.77%  '
$ws.Range('D19').Value = '5.487'
$ws.Range('E19').Value = '  -6.07%  '
$ws.Range('D20').Value = '1.003'
$ws.Range('E20').Value = '  +0.38%  '
$ws.Range('D21').Value = '210.84'
$ws.Range('E21').Value = '  -6.68%  '
$ws.Range('E22').Value = '  -1.48%  '
$ws.Range('D23').Value = '6.828'
$ws.Range('E23').Value = '  -2.30%  '
$ws.Range('E24').Value = '  +0.39%  '
$ws.Range('D25').Value = '153.83'
$ws.Range('E25').Value = '  -0.86%  '
$ws.Range('D26').Value = '7.857'
$ws.Range('E26').Value = '  -1.86%  '
$ws.Range('D27').Value = '0.1267'
$ws.Range('E27').Value = '  -2.34%  '
$ws.Range('D28').Value = '16.45'
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('D29').Value = '1.416'
$ws.Range('E29').Value = '  -4.47%  '
$ws.Range('D30').Value = '0.06212'
$ws.Range('E30').Value = '  -0.18%  '
$ws.Range('E31').Value = '  -1.42%  '
$ws.Range('D32').Value = '3.782'
$ws.Range('E32').Value = '  -1.01%  '
$ws.Range('D33').Value = '3.744'
$ws.Range('E33').Value = '  -1.41%  '
$ws.Range('E34').Value = '  -0.47%  '
$ws.Range('E35').Value = '  -5.30%  '
$ws.Range('E36').Value = '  +0.63%  '
$ws.Range('E37').Value = '  -1.19%  '
$ws.Range('D38').Value = '2.713'
$ws.Range('E38').Value = '  -0.68%  '
$ws.Range('D39').Value = '6.512'
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('D40').Value = '0.01691'
$ws.Range('E40').Value = '  -2.24%  '
$ws.Range('D41').Value = '1.148.91'
$ws.Range('E41').Value = '  -5.55%  '
$ws.Range('D42').Value = '0.8858'
$ws.Range('E42').Value = '  -2.28%  '
$ws.Range('E43').Value = '  +0.35%  '
$ws.Range('D44').Value = '99.97'
$ws.Range('E44').Value = '  -0.20%  '
$ws.Range('D45').Value = '1.951.77'
$ws.Range('E45').Value = '  -1.97%  '
$ws.Range('D46').Value = '60.58'
$ws.Range('E46').Value = '  -3.41%  '
$ws.Range('D47').Value = '''0.00000000112'
$ws.Range('E47').Value = '  -4.49%  '
$ws.Range('D48').Value = '1.592'
$ws.Range('E48').Value = '  +0.13%  '
$ws.Range('D49').Value = '8.377'
$ws.Range('E49').Value = '  -1.21%  '
$ws.Range('D50').Value = '''0.05470'
$ws.Range('E50').Value = '  -1.10%  '
$ws.Range('D51').Value = '0.4476'
$ws.Range('E51').Value = '  -1.54%  '
